$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# NOTE: shared-string indices are assigned by this runtime in the order the
# *first* occurrence of each unique string is written. To reproduce the
# exact shared-string table from the target workbook we write the brand-new
# (never-seen-before) strings in the same order they were added by the
# original author, interleaving the reused/common strings (US-09, API,
# High, Major, Critical, Pass (Postman), Pass, Robert Norwood) wherever
# convenient since those already exist in the table and will simply be
# deduplicated back to their existing index.

# --- Row 58 ---------------------------------------------------------------
$ws.Range("A58").Value = "TC-057-API-17"
$ws.Range("A59").Value = "TC-058-API-18"

$ws.Range("B58").Value = "US-09"

$ws.Range("C58").Value = "GET /api/auth/me returns 401 when not authenticated"

$ws.Range("D58").Value = "API"
$ws.Range("E58").Value = "High"
$ws.Range("F58").Value = "Major"

$ws.Range("G58").Value = "Server running GET /api/auth/me route exists; no token cookie present (logged out)"

$ws.Range("H58").Value = "In Postman, ensure cookie jar does not include token for the API domain (or start a new session).`nSend GET /api/auth/me."

$ws.Range("J58").Value = "Returns 401 Unauthorized with JSON { success:false, message:'Not authorized.' } (or your exact middleware message)"

$ws.Range("K58").Value = "Pass (Postman)"
$ws.Range("L58").Value = "Pass"
$ws.Range("M58").Value = "Robert Norwood"

$ws.Range("N58").Value = 46040

$ws.Range("O58").Value = "Manual API test via Postman. Confirms requireAuth blocks unauthenticated requests and prevents controller execution."

# --- Row 59 ---------------------------------------------------------------
$ws.Range("O59").Value = "Manual API test via Postman. Confirms session persists via HttpOnly cookie and req.user is populated from verified JWT."

$ws.Range("C59").Value = "GET /api/auth/me returns authenticated user context when logged in"

$ws.Range("B59").Value = "US-09"
$ws.Range("D59").Value = "API"
$ws.Range("E59").Value = "High"
$ws.Range("F59").Value = "Critical"

$ws.Range("G59").Value = "Admin user exists; login works; token cookie present after login"

$ws.Range("H59").Value = "Send POST /api/auth/login with valid admin credentials.`nConfirm login response is 200 and cookie token is set.`nSend GET /api/auth/me (same Postman session/cookie jar)."

$ws.Range("I59").Value = "Login body: { `"email`": `"<ADMIN_USER>`", `"password`": `"<ADMIN_PASSWORD>`" }"

$ws.Range("J59").Value = "GET /api/auth/me returns 200 OK with JSON:`nsuccess: true`ndata: { id: <string>, role: `"Admin`" }"

$ws.Range("K59").Value = "Pass (Postman)"
$ws.Range("L59").Value = "Pass"
$ws.Range("M59").Value = "Robert Norwood"

$ws.Range("N59").Value = 46040

# --- Row heights ------------------------------------------------------
$ws.Rows.Item(58).RowHeight = 75
$ws.Rows.Item(59).RowHeight = 120

# --- View state: reselect the active cell as it ended up after the edit --
$ws.Range("N64").Select()
